$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The header row ("Question" / "Answer") was mistakenly duplicated as real
# data lower in the sheet and the sheet never needed a styled header row in
# the first place. Remove the header row entirely so the data starts
# immediately with the first real question, shifting every row up by one.
$ws.Rows.Item(1).Select() | Out-Null
$ws.Rows.Item(1).Delete() | Out-Null

# Reset the view back to the top of the sheet, selecting the row that now
# occupies the first position (mirrors what Excel leaves selected after a
# whole-row delete).
$ws.Application.Goto($ws.Range("A1")) | Out-Null
$ws.Rows.Item(1).Select() | Out-Null
